$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1293.6666
$ws.Range("I17").Value = 1200
$ws.Range("K17").Value = 3600
$ws.Range("M17").Value = -3432
$ws.Range("H87").Value = 37900
$ws.Range("J87").Value = 37900
$ws.Range("L87").Value = 37900
$ws.Range("N87").Value = -40396
$ws.Range("H90").Value = 37900
$ws.Range("J90").Value = 37900
$ws.Range("L90").Value = 113700
$ws.Range("N90").Value = -126180
$ws.Range("H113").Value = 30568922
$ws.Range("J113").Value = 35731868
$ws.Range("L113").Value = 35731868
$ws.Range("N113").Value = -35738376

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1732.3914
$ws.Range("J2").Value = 1934.8
$ws.Range("L2").Value = 1934.8
$ws.Range("N2").Value = -2160.8
$ws.Range("H32").Value = 5414848
$ws.Range("I32").Value = 6256746
$ws.Range("K32").Value = 6256746
$ws.Range("M32").Value = -6256459
$ws.Range("H61").Value = 4395.7334
$ws.Range("I61").Value = 1252.9667
$ws.Range("K61").Value = 1252.9667
$ws.Range("M61").Value = -1040.9667
$ws.Range("H102").Value = 2753.5173
$ws.Range("I102").Value = 2332.8235
$ws.Range("K102").Value = 2332.8235
$ws.Range("M102").Value = -710.8235
$ws.Range("H116").Value = 1732.3914
$ws.Range("J116").Value = 1934.8
$ws.Range("L116").Value = 1934.8
$ws.Range("N116").Value = -6522.8
$ws.Range("H122").Value = 6189.643
$ws.Range("J122").Value = 6328.846
$ws.Range("L122").Value = 18986.538
$ws.Range("N122").Value = -23886.538
$ws.Range("H132").Value = 5471.9077
$ws.Range("I132").Value = 4151.2197
$ws.Range("J132").Value = 7728.0835
$ws.Range("K132").Value = 12453.6591
$ws.Range("L132").Value = 23184.2505
$ws.Range("M132").Value = -9923.659099999999
$ws.Range("N132").Value = -28244.2505
$ws.Range("H136").Value = 4395.7334
$ws.Range("I136").Value = 1252.9667
$ws.Range("K136").Value = 3758.9001
$ws.Range("M136").Value = -1208.9001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1732.3914
$ws.Range("J3").Value = 1934.8
$ws.Range("L3").Value = 1934.8
$ws.Range("N3").Value = -2162.8
$ws.Range("H94").Value = 1343.3636
$ws.Range("I94").Value = 802.5185
$ws.Range("J94").Value = 3777.1667
$ws.Range("K94").Value = 802.5185
$ws.Range("L94").Value = 3777.1667
$ws.Range("M94").Value = -351.5185
$ws.Range("N94").Value = -4679.1667
$ws.Range("H134").Value = 3650.9841
$ws.Range("I134").Value = 1256.9546
$ws.Range("K134").Value = 3770.8638
$ws.Range("M134").Value = -1235.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7100.3477
$ws.Range("I31").Value = 2040.4117
$ws.Range("J31").Value = 10066.518
$ws.Range("K31").Value = 2040.4117
$ws.Range("L31").Value = 10066.518
$ws.Range("M31").Value = -1745.4117
$ws.Range("N31").Value = -10656.518
$ws.Range("H34").Value = 7100.3477
$ws.Range("I34").Value = 2040.4117
$ws.Range("J34").Value = 10066.518
$ws.Range("K34").Value = 2040.4117
$ws.Range("L34").Value = 10066.518
$ws.Range("M34").Value = -1838.4117
$ws.Range("N34").Value = -10470.518
$ws.Range("H59").Value = 74962.664
$ws.Range("J59").Value = 74962.664
$ws.Range("L59").Value = 74962.664
$ws.Range("N59").Value = -77252.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2134.125
$ws.Range("I5").Value = 898.75
$ws.Range("J5").Value = 3369.5
$ws.Range("K5").Value = 2696.25
$ws.Range("L5").Value = 10108.5
$ws.Range("M5").Value = -2584.25
$ws.Range("N5").Value = -10332.5
$ws.Range("H69").Value = 7937.5
$ws.Range("I69").Value = 7500
$ws.Range("K69").Value = 22500
$ws.Range("M69").Value = -21689
$ws.Range("H72").Value = 7937.5
$ws.Range("I72").Value = 7500
$ws.Range("K72").Value = 67500
$ws.Range("M72").Value = -63444
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 9000
$ws.Range("M74").Value = -7939
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 27000
$ws.Range("M77").Value = -21696
$ws.Range("H113").Value = 5574.5835
$ws.Range("J113").Value = 9156.429
$ws.Range("L113").Value = 27469.287
$ws.Range("N113").Value = -31809.287
$ws.Range("H135").Value = 2134.125
$ws.Range("I135").Value = 898.75
$ws.Range("J135").Value = 3369.5
$ws.Range("K135").Value = 8088.75
$ws.Range("L135").Value = 30325.5
$ws.Range("M135").Value = -5553.75
$ws.Range("N135").Value = -35395.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 75000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 75000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -75518
$ws.Range("H58").Value = 53604.445
$ws.Range("I58").Value = 2513.3333
$ws.Range("K58").Value = 2513.3333
$ws.Range("M58").Value = -2236.3333
$ws.Range("H122").Value = 29010.021
$ws.Range("I122").Value = 35792.594
$ws.Range("K122").Value = 107377.782
$ws.Range("M122").Value = -104927.782
$ws.Range("H132").Value = 3389.6943
$ws.Range("I132").Value = 1629.037
$ws.Range("J132").Value = 8671.666999999999
$ws.Range("K132").Value = 4887.111
$ws.Range("L132").Value = 26015.001
$ws.Range("M132").Value = -2357.111
$ws.Range("N132").Value = -31075.001
$ws.Range("H134").Value = 94663
$ws.Range("J134").Value = 94663
$ws.Range("L134").Value = 283989
$ws.Range("N134").Value = -289059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 34125
$ws.Range("J2").Value = 31142.857
$ws.Range("L2").Value = 31142.857
$ws.Range("N2").Value = -31366.857
$ws.Range("H16").Value = 1062.8572
$ws.Range("I16").Value = 1090.1666
$ws.Range("J16").Value = 899
$ws.Range("K16").Value = 1090.1666
$ws.Range("L16").Value = 899
$ws.Range("M16").Value = -920.1666
$ws.Range("N16").Value = -1239
$ws.Range("H22").Value = 1617.8889
$ws.Range("I22").Value = 334
$ws.Range("J22").Value = 3222.75
$ws.Range("K22").Value = 334
$ws.Range("L22").Value = 3222.75
$ws.Range("M22").Value = -39
$ws.Range("N22").Value = -3812.75
$ws.Range("H27").Value = 1617.8889
$ws.Range("I27").Value = 334
$ws.Range("J27").Value = 3222.75
$ws.Range("K27").Value = 334
$ws.Range("L27").Value = 3222.75
$ws.Range("M27").Value = -227
$ws.Range("N27").Value = -3436.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 133338240
$ws.Range("I62").Value = 166670850
$ws.Range("K62").Value = 166670850
$ws.Range("M62").Value = -166670226
$ws.Range("H65").Value = 133338240
$ws.Range("I65").Value = 166670850
$ws.Range("K65").Value = 833354250
$ws.Range("M65").Value = -833351130
$ws.Range("H96").Value = 849.5
$ws.Range("I96").Value = 849.5
$ws.Range("K96").Value = 849.5
$ws.Range("M96").Value = 523.5
$ws.Range("H100").Value = 729.06665
$ws.Range("I100").Value = 319.25
$ws.Range("J100").Value = 1197.4286
$ws.Range("K100").Value = 638.5
$ws.Range("L100").Value = 2394.8572
$ws.Range("M100").Value = -97.5
$ws.Range("N100").Value = -3476.8572
$ws.Range("H132").Value = 5367.6875
$ws.Range("I132").Value = 5294
$ws.Range("J132").Value = 5556
$ws.Range("K132").Value = 15882
$ws.Range("L132").Value = 16668
$ws.Range("M132").Value = -13352
$ws.Range("N132").Value = -21728
$ws.Range("H136").Value = 23259640
$ws.Range("I136").Value = 41667950
$ws.Range("J136").Value = 7037.0527
$ws.Range("K136").Value = 125003850
$ws.Range("L136").Value = 21111.1581
$ws.Range("M136").Value = -125001300
$ws.Range("N136").Value = -26211.1581
